$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 31   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/12/2024  Through  8/18/2024"

# --- C14: text placeholder "0" -> numeric 1 (style should follow #,##0 numeric format) ---
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("C14").Value = 1

# --- Numeric value updates ---
# Row 14
$ws.Range("F14").Value = 3
$ws.Range("H14").Value = 50
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = 0
# Row 16
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 57
$ws.Range("J16").Value = 51
$ws.Range("K16").Value = 11.764705882352
$ws.Range("L16").Value = 42.5
# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 45.454545454545
$ws.Range("I17").Value = 134
$ws.Range("J17").Value = 143
$ws.Range("K17").Value = -6.293706293706
$ws.Range("L17").Value = -13.548387096774
# Row 18
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 33.333333333333
$ws.Range("L18").Value = -1.754385964912
# Row 19
$ws.Range("C19").Value = 8
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = -16.666666666666
$ws.Range("I19").Value = 218
$ws.Range("J19").Value = 255
$ws.Range("K19").Value = -14.509803921568
$ws.Range("L19").Value = -22.142857142857
# Row 20
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 85.714285714285
$ws.Range("I20").Value = 59
$ws.Range("J20").Value = 68
$ws.Range("K20").Value = -13.235294117647
$ws.Range("L20").Value = 3.508771929824
# Row 21
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -15
$ws.Range("F21").Value = 78
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = -1.265822784810
$ws.Range("I21").Value = 537
$ws.Range("J21").Value = 568
$ws.Range("K21").Value = -5.457746478873
$ws.Range("L21").Value = -9.747899159663
# Row 24
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -12.121212121212
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 126
$ws.Range("H24").Value = -16.666666666666
$ws.Range("I24").Value = 907
$ws.Range("J24").Value = 892
$ws.Range("K24").Value = 1.681614349775
$ws.Range("L24").Value = 3.302961275626
# Row 25
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 67
$ws.Range("H25").Value = -4.477611940298
$ws.Range("I25").Value = 614
$ws.Range("J25").Value = 521
$ws.Range("K25").Value = 17.850287907869
$ws.Range("L25").Value = 41.149425287356
# Row 26
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = -55.555555555555
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = -20.454545454545
$ws.Range("I26").Value = 376
$ws.Range("J26").Value = 347
$ws.Range("K26").Value = 8.357348703170
$ws.Range("L26").Value = 14.634146341463
# Row 28
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("J28").Value = 39
$ws.Range("K28").Value = -5.128205128205
# Row 29
$ws.Range("F29").Value = 2
# Row 30
$ws.Range("F30").Value = 1

# --- Cells switching from numeric to text placeholders ---
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("C23").Value = "'0"
$ws.Range("C27").Value = "'0"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "'***.*"
$ws.Range("C28").Value = "'0"
